$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at the top of the data (rows 2 and 3), pushing
# everything else down by two rows.
$ws.Range("A2:A3").EntireRow.Insert()

# Row 2: new "Authentication Testing" suite header.
$ws.Range("A2").Value = "Authentication Testing"
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""

# Copy the look of an existing suite-header row (still A4:C4, the old
# "User Testing" header) onto the new header row.
$ws.Range("A4:C4").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)

# Row 3: new "Authentication Testing - Log-in" test-case header row.
$ws.Range("A3").Value = "  Authentication Testing - Log-in"
$ws.Range("B3").Value = "AUTH-0001"
$ws.Range("C3").Value = "Verify that user can log-in successfully. "

# Base the new row's look on the existing bold sub-header styles used
# elsewhere in the sheet (column A / B / C each have their own look),
# then add the same white/"theme background 1" fill used by the author.
# (Rows shifted down by 2 after the insert above: old row 39 -> 41,
# old row 46 -> 48.)
$ws.Range("A48").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("B41").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("C48").Copy()
$ws.Range("C3").PasteSpecial(-4122)

$row3 = $ws.Range("A3:C3")
$row3.Interior.PatternColor = 13888217
$row3.Interior.ThemeColor = 2
$row3.Interior.Pattern = 1

# Update view: selection moves to D5, and the previous "frozen-in-place"
# top-left scroll position is cleared.
$ws.Range("D5").Select()
